$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 17:14:16'
$ws.Range('A3').Value = 'Total filas: 291'
$ws.Range('A16').Value = '05:44:02'
$ws.Range('C16').Value = '17X38_ROMERO'
$ws.Range('D16').Value = 56
$ws.Range('A17').Value = '06:38:54'
$ws.Range('C17').Value = '16_SANTA ANA'
$ws.Range('D17').Value = 2
$ws.Range('A41').Value = '07:52:32'
$ws.Range('C41').Value = '17_ROMERO'
$ws.Range('D41').Value = 8
$ws.Range('A42').Value = '06:38:54'
$ws.Range('C42').Value = '16_SANTA ANA'
$ws.Range('D42').Value = 82
$ws.Range('C49').Value = '11_ETCHEVERRY'
$ws.Range('C50').Value = '15_ABASTO'
$ws.Range('C51').Value = '11_ETCHEVERRY'
$ws.Range('C52').Value = '15_ABASTO'
$ws.Range('C64').Value = '215B_EL PATO'
$ws.Range('C65').Value = '23_HERNANDEZ'
$ws.Range('A75').Value = '08:52:33'
$ws.Range('C75').Value = '14_ABASTO'
$ws.Range('D75').Value = 25
$ws.Range('A76').Value = '08:30:14'
$ws.Range('C76').Value = '27_EL RETIRO'
$ws.Range('D76').Value = 47
$ws.Range('A77').Value = '08:40:59'
$ws.Range('C77').Value = '15X38_ABASTO'
$ws.Range('D77').Value = 37
$ws.Range('A85').Value = '08:30:14'
$ws.Range('C85').Value = '11_ETCHEVERRY'
$ws.Range('D85').Value = 72
$ws.Range('A86').Value = '08:40:59'
$ws.Range('C86').Value = '16_P MOR-SANTA ANA'
$ws.Range('D86').Value = 62
$ws.Range('C102').Value = '14_ABASTO'
$ws.Range('C103').Value = '15_ABASTO'
$ws.Range('A112').Value = '09:23:52'
$ws.Range('C112').Value = '10_OLMOS'
$ws.Range('D112').Value = 98
$ws.Range('A113').Value = '10:56:01'
$ws.Range('C113').Value = '81_EL PELIGRO'
$ws.Range('D113').Value = 5
$ws.Range('C149').Value = '23_HERNANDEZ'
$ws.Range('C150').Value = '14_ABASTO'
$ws.Range('A153').Value = '12:33:54'
$ws.Range('C153').Value = '15_ABASTO'
$ws.Range('D153').Value = 1
$ws.Range('A154').Value = '10:56:01'
$ws.Range('C154').Value = '27_EL RETIRO'
$ws.Range('D154').Value = 98
$ws.Range('A174').Value = '13:14:41'
$ws.Range('C174').Value = '14_ABASTO'
$ws.Range('D174').Value = 18
$ws.Range('A175').Value = '12:33:54'
$ws.Range('C175').Value = '215A_EL PATO'
$ws.Range('D175').Value = 59
$ws.Range('A186').Value = '13:57:31'
$ws.Range('C186').Value = '10_OLMOS'
$ws.Range('D186').Value = 5
$ws.Range('A187').Value = '13:14:41'
$ws.Range('C187').Value = '23_HERNANDEZ'
$ws.Range('D187').Value = 48
$ws.Range('A188').Value = '13:43:25'
$ws.Range('C188').Value = '16_SANTA ANA'
$ws.Range('D188').Value = 19
$ws.Range('C201').Value = '215C_EL PATO'
$ws.Range('C202').Value = '14X44_ABASTO'
$ws.Range('A223').Value = '15:17:56'
$ws.Range('C223').Value = '10_OLMOS'
$ws.Range('D223').Value = 35
$ws.Range('A224').Value = '14:53:58'
$ws.Range('C224').Value = '16_P MOR-SANTA ANA'
$ws.Range('D224').Value = 59
$ws.Range('A225').Value = '14:53:58'
$ws.Range('C225').Value = '15X38_ABASTO'
$ws.Range('D225').Value = 59
$ws.Range('A226').Value = '14:33:43'
$ws.Range('C226').Value = '27_EL RETIRO'
$ws.Range('D226').Value = 79
$ws.Range('A229').Value = '15:47:47'
$ws.Range('C229').Value = '15X38_ABASTO'
$ws.Range('D229').Value = 6
$ws.Range('A230').Value = '14:47:05'
$ws.Range('C230').Value = '27_EL RETIRO'
$ws.Range('D230').Value = 66
$ws.Range('A262').Value = '17:14:15'
$ws.Range('D262').Value = 7
$ws.Range('A265').Value = '17:14:15'
$ws.Range('D265').Value = 20
$ws.Range('A267').Value = '17:14:15'
$ws.Range('D267').Value = 23
$ws.Range('A269').Value = '17:14:15'
$ws.Range('D269').Value = 25
$ws.Range('A272').Value = '17:14:15'
$ws.Range('D272').Value = 32
$ws.Range('A275').Value = '17:14:15'
$ws.Range('D275').Value = 37
$ws.Range('A278').Value = '17:14:15'
$ws.Range('D278').Value = 39
$ws.Range('A279').Value = '17:14:15'
$ws.Range('D279').Value = 44
$ws.Range('A281').Value = '17:14:15'
$ws.Range('D281').Value = 52
$ws.Range('A282').Value = '17:14:15'
$ws.Range('D282').Value = 56
$ws.Range('A283').Value = '17:14:15'
$ws.Range('D283').Value = 56
$ws.Range('A286').Value = '17:14:15'
$ws.Range('D286').Value = 68
$ws.Range('A287').Value = '17:14:15'
$ws.Range('D287').Value = 71
$ws.Range('A289').Value = '17:14:15'
$ws.Range('D289').Value = 76
$ws.Range('A290').Value = '17:14:15'
$ws.Range('D290').Value = 82
$ws.Range('A291').Value = '17:14:15'
$ws.Range('B291').Value = '18:36'
$ws.Range('C291').Value = '23_HERNANDEZ'
$ws.Range('D291').Value = 82
$ws.Range('E291').Value = 'LP1912'
$ws.Range('A292').Value = '17:14:15'
$ws.Range('B292').Value = '18:41'
$ws.Range('C292').Value = '10_OLMOS'
$ws.Range('D292').Value = 87
$ws.Range('E292').Value = 'LP1912'
$ws.Range('A293').Value = '17:14:15'
$ws.Range('B293').Value = '18:45'
$ws.Range('C293').Value = '16_SANTA ANA'
$ws.Range('D293').Value = 91
$ws.Range('E293').Value = 'LP1912'
$ws.Range('A294').Value = '17:14:15'
$ws.Range('B294').Value = '18:57'
$ws.Range('C294').Value = '16_P MOR-SANTA ANA'
$ws.Range('D294').Value = 103
$ws.Range('E294').Value = 'LP1912'
$ws.Range('A295').Value = '17:14:15'
$ws.Range('B295').Value = '19:00'
$ws.Range('C295').Value = '14_ABASTO'
$ws.Range('D295').Value = 106
$ws.Range('E295').Value = 'LP1912'
$ws.Range('A296').Value = '17:14:15'
$ws.Range('B296').Value = '19:04'
$ws.Range('C296').Value = '215_EL PELIGRO'
$ws.Range('D296').Value = 110
$ws.Range('E296').Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 17:14:16'
$ws.Range('A3').Value = 'Total filas: 43'
$ws.Range('A42').Value = '17:14:15'
$ws.Range('D42').Value = 25
$ws.Range('A45').Value = '17:14:15'
$ws.Range('D45').Value = 32
$ws.Range('A47').Value = '17:14:15'
$ws.Range('D47').Value = 68
$ws.Range('A48').Value = '17:14:15'
$ws.Range('B48').Value = '19:04'
$ws.Range('C48').Value = '215_EL PELIGRO'
$ws.Range('D48').Value = 110
$ws.Range('E48').Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 17:14:16'
$ws.Range('A43').Value = '17:14:15'
$ws.Range('D43').Value = 82
